$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on all cells we touch, so numeric-looking
# strings (e.g. "29.898.72", "1.000") are preserved literally as text
# rather than being auto-converted to numbers by Excel's type inference.
$cells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12",
    "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17",
    "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22",
    "D23", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28",
    "E28", "D29", "E29", "D30", "E30", "B31", "C31", "D31", "E31", "B32",
    "C32", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36",
    "E37", "D38", "E38", "E39", "D40", "E40", "D41", "E41", "D42", "E42",
    "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "D48", "E48",
    "D49", "E49", "D50", "E50", "E51"
)
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "29.898.72"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.888.74"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "0.7734"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "242.79"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "25.66"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").Value = "0.07187"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "0.08616"
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("D12").Value = "1.953.76"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").Value = "0.7644"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "5.376"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "93.83"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "6.185"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "29.978.65"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "13.79"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "244.69"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "0.000007828"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "2.197.30"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "8.017"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "0.1648"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "9.381"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "161.97"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "18.76"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "2.034"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "1.444"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.534"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.529"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "4.103"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "0.05427"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "0.7456"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").Value = "2.781"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "0.4467"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "1.111.03"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("D43").Value = "73.34"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "6.085"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").Value = "0.8514"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "102.86"
$ws.Range("D48").Value = "1.871"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "7.630"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "2.092.55"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").Value = "  -0.78%  "
